# Updates the crypto price/volume table with refreshed values.
# Commit: Updated cryptos list on Fri Feb 23 18:07:42 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.088.46'
$ws.Range("E2").Value = '  -1.01%  '
$ws.Range("D3").Value = '2.942.44'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '374.76'
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.26'
$ws.Range("E6").Value = '  -3.02%  '
$ws.Range("E7").Value = '  -1.25%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.26'
$ws.Range("E10").Value = '  -2.75%  '
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0853'
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").Value = '3.411.84'
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.02'
$ws.Range("E14").Value = '  -2.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.54'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '2.939.54'
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '10.99'
$ws.Range("E17").Value = '  +47.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.980'
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("D19").Value = '51.083.77'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.12'
$ws.Range("E20").Value = '  -6.16%  '
$ws.Range("E21").Value = '  -4.26%  '
$ws.Range("D22").Value = '0.0₃0959'
$ws.Range("E22").Value = '  -0.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '265.33'
$ws.Range("E23").Value = '  +1.04%  '
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.14'
$ws.Range("E25").Value = '  +9.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.10'
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.57'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.62'
$ws.Range("E29").Value = '  -1.15%  '
$ws.Range("E30").Value = '  -4.55%  '
$ws.Range("E31").Value = '  -5.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.99'
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.73'
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.05'
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.44'
$ws.Range("E35").Value = '  -5.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0443'
$ws.Range("E36").Value = '  -1.45%  '
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.17'
$ws.Range("E38").Value = '  +4.09%  '
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.26'
$ws.Range("E40").Value = '  -5.20%  '
$ws.Range("E41").Value = '  -3.67%  '
$ws.Range("E42").Value = '  -3.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.05'
$ws.Range("E43").Value = '  -4.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.29'
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.03'
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.33'
$ws.Range("E46").Value = '  +2.39%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.272'
$ws.Range("E47").Value = '  -4.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.31'
$ws.Range("E48").Value = '  -3.29%  '
$ws.Range("D49").Value = '1.996.70'
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0324'
$ws.Range("E50").Value = '  -2.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.32'
$ws.Range("E51").Value = '  +2.00%  '
